# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# sheets to reflect newly scraped totals (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 168
$wsExhibit.Range("F4").Value  = 79
$wsExhibit.Range("F7").Value  = 1655
$wsExhibit.Range("F8").Value  = 9
$wsExhibit.Range("F9").Value  = 12
$wsExhibit.Range("F10").Value = 23
$wsExhibit.Range("F11").Value = 1551
$wsExhibit.Range("F12").Value = 123
$wsExhibit.Range("F13").Value = 50
$wsExhibit.Range("F14").Value = 382
$wsExhibit.Range("F20").Value = 47
$wsExhibit.Range("F21").Value = 275
$wsExhibit.Range("F22").Value = 154
$wsExhibit.Range("F23").Value = 217

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 168
$wsAll.Range("F4").Value  = 79
$wsAll.Range("F7").Value  = 1655
$wsAll.Range("F9").Value  = 9
$wsAll.Range("F10").Value = 12
$wsAll.Range("F11").Value = 23
$wsAll.Range("F12").Value = 1551
$wsAll.Range("F13").Value = 123
$wsAll.Range("F14").Value = 50
$wsAll.Range("F15").Value = 382
$wsAll.Range("F21").Value = 47
$wsAll.Range("F22").Value = 275
$wsAll.Range("F23").Value = 154
$wsAll.Range("F24").Value = 217
